$d = $word.ActiveDocument

$pairs = @(
    @("2026-01-01 Thursday", "2026-01-02 Friday"),
    @("74÷8=9, 2", "87÷8=10, 7"),
    @("79÷3=26, 1", "85÷7=12, 1"),
    @("31÷6=5, 1", "91÷6=15, 1"),
    @("40÷5=8, 0", "51÷6=8, 3"),
    @("98÷4=24, 2", "11÷7=1, 4"),
    @("48÷4=12, 0", "23÷6=3, 5"),
    @("40÷4=10, 0", "77÷5=15, 2"),
    @("66÷4=16, 2", "80÷7=11, 3"),
    @("59÷7=8, 3", "16÷2=8, 0"),
    @("95÷5=19, 0", "88÷6=14, 4"),
    @("91÷3=30, 1", "49÷2=24, 1"),
    @("98÷2=49, 0", "94÷4=23, 2"),
    @("15÷8=1, 7", "17÷7=2, 3"),
    @("36÷9=4, 0", "65÷3=21, 2"),
    @("56÷4=14, 0", "15÷7=2, 1"),
    @("31÷5=6, 1", "83÷6=13, 5"),
    @("88÷6=14, 4", "43÷3=14, 1"),
    @("15÷3=5, 0", "15÷9=1, 6"),
    @("31÷2=15, 1", "87÷7=12, 3"),
    @("29÷4=7, 1", "48÷9=5, 3"),
    @("76÷2=38, 0", "36÷5=7, 1"),
    @("61÷5=12, 1", "54÷4=13, 2"),
    @("62÷5=12, 2", "88÷3=29, 1"),
    @("98÷2=49, 0", "84÷8=10, 4"),
    @("47÷9=5, 2", "83÷6=13, 5")
)

$searchStart = 0
foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $r = $d.Range($searchStart, $d.Content.End)
    $found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 1)
    if (-not $found) {
        Write-Host "NOT FOUND:" $old
    } else {
        $searchStart = $r.End
    }
}
Write-Host "Done"
